$wb = $excel.ActiveWorkbook

# --- Break the external link to the (missing) source workbook ("SHEEt") ---
# This both removes the <externalReferences> from the workbook part and
# converts every formula that referenced it ( [1]Materiais!A:A ) on the
# "Quantidades" sheet into the plain cached value it already held.
$wb.BreakLink("SHEEt", 1) | Out-Null

$ws1 = $wb.Worksheets.Item("Materiais")
$ws2 = $wb.Worksheets.Item("Quantidades")

# --- The row-24 material code lookup is missing now - clear the leftover
# cell entirely (code 7009467 row), keeping the Qty/Stock columns intact ---
$ws2.Range("A24").ClearContents() | Out-Null

# --- Explicit portrait page setup on the Materiais sheet ---
$ws1.PageSetup.Orientation = 1

# --- Update selections / active sheet to match where editing left off ---
$ws1.Activate() | Out-Null
$ws1.Range("B1").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A24").Select() | Out-Null
